$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 117; existing rows 117-136 shift down to 118-137.
$ws.Rows(117).Insert()

# Populate the newly inserted row 117 with the new weekly price record.
$ws.Cells.Item(117, 1).Value = 7
$ws.Cells.Item(117, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(117, 3).Value = "Ñuble"
$ws.Cells.Item(117, 4).Value = 45212
$ws.Cells.Item(117, 5).Value = 16
$ws.Cells.Item(117, 6).Value = 100112013
$ws.Cells.Item(117, 7).Value = "Alcachofa"
$ws.Cells.Item(117, 8).Value = "Española"
$ws.Cells.Item(117, 9).Value = "Primera"
$ws.Cells.Item(117, 10).Value = 110
$ws.Cells.Item(117, 11).Value = 15000
$ws.Cells.Item(117, 12).Value = 15500
$ws.Cells.Item(117, 13).Value = 15227
$ws.Cells.Item(117, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(117, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(117, 16).Value = 508
$ws.Cells.Item(117, 17).Value = 30
$ws.Cells.Item(117, 18).Value = "Hortaliza"
